# Apply updated crypto price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.839.39'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '1.635.70'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '''215.30'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('D6').Value = '''0.5016'
$ws.Range('E6').Value = '  -2.74%  '
$ws.Range('D7').Value = '''1.003'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('D9').Value = '''0.06414'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '''19.57'
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('D11').Value = '''0.07679'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('D12').Value = '1.657.59'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '''4.242'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').Value = '1.861.60'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').Value = '0.0₅7926'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').Value = '''63.50'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '25.857.66'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = '''203.05'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').Value = '''4.301'
$ws.Range('E21').Value = '  -2.70%  '
$ws.Range('D22').Value = '''9.945'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('D23').Value = '''5.982'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = '''1.937'
$ws.Range('E25').Value = '  +10.20%  '
$ws.Range('D26').Value = '''141.38'
$ws.Range('E26').Value = '  -2.29%  '
$ws.Range('D27').Value = '''0.1145'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').Value = '''6.703'
$ws.Range('E29').Value = '  -3.91%  '
$ws.Range('D30').Value = '''1.239'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').Value = '''0.04967'
$ws.Range('E31').Value = '  -6.01%  '
$ws.Range('D32').Value = '''3.262'
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').Value = '''3.184'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('D34').Value = '''1.528'
$ws.Range('E34').Value = '  -2.82%  '
$ws.Range('D35').Value = '''2.350'
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('D36').Value = '1.174.22'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('D38').Value = '''2.621'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('D40').Value = '''0.01557'
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').Value = '''2.554'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = '''5.640'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '''0.8045'
$ws.Range('E44').Value = '  -4.72%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''99.50'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').Value = '1.774.25'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').Value = '0.0₈110'
$ws.Range('E47').Value = '  -3.78%  '
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '''54.83'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').Value = '''0.05043'
$ws.Range('E51').Value = '  -0.22%  '
